# The BTec header logo (image1.jpg) is renamed to image2.jpg, and the
# Pearson footer logo (image2.png) is renamed to image1.png, in both the
# first-page and default header/footer parts.
$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers($i)
        if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }

        $ftr = $sec.Footers($i)
        if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.AlternativeText -like "*PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
